$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data rows 2:28 (A:I) before rewriting in sorted order
$ws.Range("A2:I28").ClearContents()

$ws.Range("A2").Value = "FALX99"
$ws.Range("B2").Value = "Course"
$ws.Range("C2").Value = "EDUC"
$ws.Range("D2").Value = "FAL X99 - Foundations of Academic Literacy"
$ws.Range("E2").Value = "None"
$ws.Range("F2").Value = "None"
$ws.Range("G2").Value = "EDUC"
$ws.Range("H2").Value = "EDUCATION"

$ws.Range("A3").Value = "GEOG221"
$ws.Range("B3").Value = "Course"
$ws.Range("C3").Value = "GEOG"
$ws.Range("D3").Value = "GEOG 221 - Economic Worlds"
$ws.Range("E3").Value = "GEOG100"
$ws.Range("F3").Value = "None"
$ws.Range("G3").Value = "ENV"
$ws.Range("H3").Value = "GEOGRAPH"
$ws.Range("I3").Value = "REQ-GEOG 100."

$ws.Range("A4").Value = "LBST100"
$ws.Range("B4").Value = "Course"
$ws.Range("C4").Value = "LBST"
$ws.Range("D4").Value = "LBST 100 - Equality and Inequality at Work"
$ws.Range("E4").Value = "None"
$ws.Range("F4").Value = "None"
$ws.Range("G4").Value = "ARTS"
$ws.Range("H4").Value = "LABOUR STU"

$ws.Range("A5").Value = "LBST101"
$ws.Range("B5").Value = "Course"
$ws.Range("C5").Value = "LBST"
$ws.Range("D5").Value = "LBST 101 - Work and Worker's Rights: Introducing Labour Studies"
$ws.Range("E5").Value = "None"
$ws.Range("F5").Value = "None"
$ws.Range("G5").Value = "ARTS"
$ws.Range("H5").Value = "LABOUR STU"

$ws.Range("A6").Value = "LBST201"
$ws.Range("B6").Value = "Course"
$ws.Range("C6").Value = "LBST"
$ws.Range("D6").Value = "LBST 201 - Workers in the Global Economy: Globalization, Labour and Uneven Development"
$ws.Range("E6").Value = "None"
$ws.Range("F6").Value = "None"
$ws.Range("G6").Value = "ARTS"
$ws.Range("H6").Value = "LABOUR STU"

$ws.Range("A7").Value = "LBST202"
$ws.Range("B7").Value = "Course"
$ws.Range("C7").Value = "LBST"
$ws.Range("D7").Value = "LBST 202 - Labour Research for Social Change: Methods and Approaches"
$ws.Range("E7").Value = "FANX99,LBST100,LBST101"
$ws.Range("F7").Value = "None"
$ws.Range("G7").Value = "ARTS"
$ws.Range("H7").Value = "LABOUR STU"
$ws.Range("I7").Value = "Prerequisite: LBST 100 or LBST 101 or with permission of instructor.   Quantitative."

$ws.Range("A8").Value = "LBST203"
$ws.Range("B8").Value = "Course"
$ws.Range("C8").Value = "LBST"
$ws.Range("D8").Value = "LBST 203 - Work and Health"
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = "None"
$ws.Range("G8").Value = "ARTS"
$ws.Range("H8").Value = "LABOUR STU"

$ws.Range("A9").Value = "LBST230"
$ws.Range("B9").Value = "Course"
$ws.Range("C9").Value = "LBST"
$ws.Range("D9").Value = "LBST 230 - Special Topics in Labour Studies"
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = "None"
$ws.Range("G9").Value = "ARTS"
$ws.Range("H9").Value = "LABOUR STU"

$ws.Range("A10").Value = "LBST301W"
$ws.Range("B10").Value = "Course"
$ws.Range("C10").Value = "LBST"
$ws.Range("D10").Value = "LBST 301W - Labour, Social Media, and the News"
$ws.Range("E10").Value = "FALX99,LBST100,LBST101"
$ws.Range("F10").Value = "None"
$ws.Range("G10").Value = "ARTS"
$ws.Range("H10").Value = "LABOUR STU"
$ws.Range("I10").Value = "REQ- LBST 100 or LBST 101 or on permission of instructor.  Students with credit for LBST 301 may not take this course for further credit. Writing."

$ws.Range("A11").Value = "LBST305"
$ws.Range("B11").Value = "Course"
$ws.Range("C11").Value = "LBST"
$ws.Range("D11").Value = "LBST 305 - Gendering Economy: Paid and Unpaid Labour"
$ws.Range("E11").Value = "LBST100"
$ws.Range("F11").Value = "None"
$ws.Range("G11").Value = "ARTS"
$ws.Range("H11").Value = "LABOUR STU"
$ws.Range("I11").Value = "Prerequisite: 30 units including LBST 100 or three units in GSWS or WS or GDST."

$ws.Range("A12").Value = "LBST306"
$ws.Range("B12").Value = "Course"
$ws.Range("C12").Value = "LBST"
$ws.Range("D12").Value = "LBST 306 - The Political Economy of Labour Markets: Critical and Radical Approaches"
$ws.Range("E12").Value = "None"
$ws.Range("F12").Value = "None"
$ws.Range("G12").Value = "ARTS"
$ws.Range("H12").Value = "LABOUR STU"

$ws.Range("A13").Value = "LBST307"
$ws.Range("B13").Value = "Course"
$ws.Range("C13").Value = "LBST"
$ws.Range("D13").Value = "LBST 307 - Unfree Labour and Modern Slavery: Understanding Coercion and Exploitation in Contemporary"
$ws.Range("E13").Value = "None"
$ws.Range("F13").Value = "None"
$ws.Range("G13").Value = "ARTS"
$ws.Range("H13").Value = "LABOUR STU"

$ws.Range("A14").Value = "LBST308"
$ws.Range("B14").Value = "Course"
$ws.Range("C14").Value = "LBST"
$ws.Range("D14").Value = "LBST 308 - The Labour Process: Technological Change and the Future of Work"
$ws.Range("E14").Value = "None"
$ws.Range("F14").Value = "None"
$ws.Range("G14").Value = "ARTS"
$ws.Range("H14").Value = "LABOUR STU"

$ws.Range("A15").Value = "LBST309"
$ws.Range("B15").Value = "Course"
$ws.Range("C15").Value = "LBST"
$ws.Range("D15").Value = "LBST 309 - Labour and Collective Bargaining"
$ws.Range("E15").Value = "None"
$ws.Range("F15").Value = "None"
$ws.Range("G15").Value = "ARTS"
$ws.Range("H15").Value = "LABOUR STU"

$ws.Range("A16").Value = "LBST310"
$ws.Range("B16").Value = "Course"
$ws.Range("C16").Value = "LBST"
$ws.Range("D16").Value = "LBST 310 - The Politics of Labour"
$ws.Range("E16").Value = "None"
$ws.Range("F16").Value = "None"
$ws.Range("G16").Value = "ARTS"
$ws.Range("H16").Value = "LABOUR STU"
$ws.Range("I16").Value = "REQ-30 units. Strongly Recommended: LBST 101."

$ws.Range("A17").Value = "LBST311"
$ws.Range("B17").Value = "Course"
$ws.Range("C17").Value = "LBST"
$ws.Range("D17").Value = "LBST 311 - Labour and the Environment"
$ws.Range("E17").Value = "None"
$ws.Range("F17").Value = "None"
$ws.Range("G17").Value = "ARTS"
$ws.Range("H17").Value = "LABOUR STU"

$ws.Range("A18").Value = "LBST312"
$ws.Range("B18").Value = "Course"
$ws.Range("C18").Value = "LBST"
$ws.Range("D18").Value = "LBST 312 - Migration, Migrants, and Work: A Global Perspective"
$ws.Range("E18").Value = "None"
$ws.Range("F18").Value = "None"
$ws.Range("G18").Value = "ARTS"
$ws.Range("H18").Value = "LABOUR STU"

$ws.Range("A19").Value = "LBST313"
$ws.Range("B19").Value = "Course"
$ws.Range("C19").Value = "LBST"
$ws.Range("D19").Value = "LBST 313 - Introduction to Canadian Labour Law"
$ws.Range("E19").Value = "None"
$ws.Range("F19").Value = "None"
$ws.Range("G19").Value = "ARTS"
$ws.Range("H19").Value = "LABOUR STU"

$ws.Range("A20").Value = "LBST320"
$ws.Range("B20").Value = "Course"
$ws.Range("C20").Value = "LBST"
$ws.Range("D20").Value = "LBST 320 - Labour and Popular Culture: Class, Politics, and Pop Culture"
$ws.Range("E20").Value = "None"
$ws.Range("F20").Value = "None"
$ws.Range("G20").Value = "ARTS"
$ws.Range("H20").Value = "LABOUR STU"
$ws.Range("I20").Value = "REQ-30 units.  Students with credit for LBST 330 under the title `"Labour and Film`" may not take this course for further credit."

$ws.Range("A21").Value = "LBST328"
$ws.Range("B21").Value = "Course"
$ws.Range("C21").Value = "LBST"
$ws.Range("D21").Value = "LBST 328 - Labour Geographies"
$ws.Range("E21").Value = "GEOG221,LBST101"
$ws.Range("F21").Value = "None"
$ws.Range("G21").Value = "ARTS"
$ws.Range("H21").Value = "LABOUR STU"
$ws.Range("I21").Value = "REQ-60 units; LBST 101 or GEOG 221.  Students with credit for GEOG 328 may not take this course for further credit."

$ws.Range("A22").Value = "LBST330"
$ws.Range("B22").Value = "Course"
$ws.Range("C22").Value = "LBST"
$ws.Range("D22").Value = "LBST 330 - Selected Topics in Labour Studies"
$ws.Range("E22").Value = "None"
$ws.Range("F22").Value = "None"
$ws.Range("G22").Value = "ARTS"
$ws.Range("H22").Value = "LABOUR STU"

$ws.Range("A23").Value = "LBST331"
$ws.Range("B23").Value = "Course"
$ws.Range("C23").Value = "LBST"
$ws.Range("D23").Value = "LBST 331 - Selected Topics in Labour Studies"
$ws.Range("E23").Value = "None"
$ws.Range("F23").Value = "None"
$ws.Range("G23").Value = "ARTS"
$ws.Range("H23").Value = "LABOUR STU"

$ws.Range("A24").Value = "LBST401"
$ws.Range("B24").Value = "Course"
$ws.Range("C24").Value = "LBST"
$ws.Range("D24").Value = "LBST 401 - How to Make Change: Community-Labour Organizing and Action"
$ws.Range("E24").Value = "None"
$ws.Range("F24").Value = "None"
$ws.Range("G24").Value = "ARTS"
$ws.Range("H24").Value = "LABOUR STU"

$ws.Range("A25").Value = "LBST431"
$ws.Range("B25").Value = "Course"
$ws.Range("C25").Value = "LBST"
$ws.Range("D25").Value = "LBST 431 - Selected Topics in Labour Studies"
$ws.Range("E25").Value = "None"
$ws.Range("F25").Value = "None"
$ws.Range("G25").Value = "ARTS"
$ws.Range("H25").Value = "LABOUR STU"

$ws.Range("A26").Value = "LBST490"
$ws.Range("B26").Value = "Course"
$ws.Range("C26").Value = "LBST"
$ws.Range("D26").Value = "LBST 490 - Directed Readings in Labour Studies"
$ws.Range("E26").Value = "None"
$ws.Range("F26").Value = "None"
$ws.Range("G26").Value = "ARTS"
$ws.Range("H26").Value = "LABOUR STU"

$ws.Range("A27").Value = "FANX99"
$ws.Range("B27").Value = "Course"
$ws.Range("C27").Value = "MATH"
$ws.Range("D27").Value = "FAN X99 - Foundations of Analytical and Quantitative Reasoning"
$ws.Range("E27").Value = "None"
$ws.Range("F27").Value = "None"
$ws.Range("G27").Value = "SCI"
$ws.Range("H27").Value = "MATHEMATIC"

$ws.Range("A28").Value = "FANX99"
$ws.Range("B28").Value = "Course"
$ws.Range("C28").Value = "MATH"
$ws.Range("D28").Value = "FAN X99 - Foundations of Analytical and Quantitative Reasoning"
$ws.Range("E28").Value = "None"
$ws.Range("F28").Value = "None"
$ws.Range("G28").Value = "SCI"
$ws.Range("H28").Value = "MATHEMATIC"
$ws.Range("I28").Value = "REQ-students who have taken, have received transfer credit for, or are currently taking MATH 150, 151, 154 or 157 may not take FAN X99 for credit without the permission from the Department of Mathematics."


# Update selection to match the saved view state
$ws.Range("A2:I28").Select()
